$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 1969.8
$ws.Range("D2").Value = 325.4
$ws.Range("E2").Value = 17576.3
$ws.Range("F2").Value = 4617.1
$ws.Range("G2").Value = 8.92
$ws.Range("H2").Value = 0.75
$ws.Range("J2").Value = 1245.1
$ws.Range("K2").Value = 327.1
$ws.Range("L2").Value = "17576.3 (±4617.1)"
$ws.Range("M2").Value = "8.92 (±0.75)"
$ws.Range("N2").Value = "1245.1 (±327.1)"

# Row 3
$ws.Range("C3").Value = 4052.400000000001
$ws.Range("D3").Value = 369.8000000000001
$ws.Range("E3").Value = 44097.4
$ws.Range("F3").Value = 5709.2
$ws.Range("G3").Value = 10.88
$ws.Range("H3").Value = 0.38
$ws.Range("J3").Value = 3131.7
$ws.Range("K3").Value = 405.5
$ws.Range("L3").Value = "44097.4 (±5709.2)"
$ws.Range("M3").Value = "10.88 (±0.38)"
$ws.Range("N3").Value = "3131.7 (±405.5)"

# Row 4
$ws.Range("C4").Value = 6020.8
$ws.Range("D4").Value = 560.6999999999999
$ws.Range("E4").Value = 61632.1
$ws.Range("F4").Value = 8755.100000000002
$ws.Range("G4").Value = 10.24
$ws.Range("H4").Value = 0.45
$ws.Range("J4").Value = 2185.8
$ws.Range("K4").Value = 310.5
$ws.Range("L4").Value = "61632.1 (±8755.1)"
$ws.Range("M4").Value = "10.24 (±0.45)"
$ws.Range("N4").Value = "2185.8 (±310.5)"

# Row 11
$ws.Range("C11").Value = 9928.0
$ws.Range("E11").Value = 74966.0
$ws.Range("J11").Value = 1284.2
$ws.Range("L11").Value = "74966.0 (±16832.1)"
$ws.Range("N11").Value = "1284.2 (±288.3)"

# Row 12
$ws.Range("C12").Value = 9220.599999999999
$ws.Range("E12").Value = 82101.50000000001
$ws.Range("J12").Value = 1446.1
$ws.Range("L12").Value = "82101.5 (±19380.7)"
$ws.Range("N12").Value = "1446.1 (±341.4)"

# Row 13
$ws.Range("C13").Value = 19108.4
$ws.Range("E13").Value = 158352.6
$ws.Range("J13").Value = 1375.2
$ws.Range("L13").Value = "158352.6 (±34108.9)"
$ws.Range("N13").Value = "1375.2 (±296.2)"

# Row 14
$ws.Range("C14").Value = 7308.4
$ws.Range("E14").Value = 90977.4
$ws.Range("J14").Value = 2542.7
$ws.Range("L14").Value = "90977.4 (±15022.6)"
$ws.Range("N14").Value = "2542.7 (±419.9)"

# Row 16
$ws.Range("C16").Value = 17593.2
$ws.Range("E16").Value = 221671.0
$ws.Range("J16").Value = 3191.2
$ws.Range("L16").Value = "221671.0 (±31090.5)"
$ws.Range("N16").Value = "3191.2 (±447.6)"

# Row 23
$ws.Range("C23").Value = 8915.0
$ws.Range("E23").Value = 84348.09999999999
$ws.Range("J23").Value = 1556.1
$ws.Range("L23").Value = "84348.1 (±20195.3)"
$ws.Range("N23").Value = "1556.1 (±372.6)"

# Row 24
$ws.Range("C24").Value = 10774.2
$ws.Range("E24").Value = 108327.9
$ws.Range("H24").Value = 0.44
$ws.Range("J24").Value = 2055.3
$ws.Range("L24").Value = "108327.9 (±21620.1)"
$ws.Range("M24").Value = "10.05 (±0.44)"
$ws.Range("N24").Value = "2055.3 (±410.2)"

# Row 25
$ws.Range("C25").Value = 19679.8
$ws.Range("E25").Value = 198531.1
$ws.Range("J25").Value = 1857.0
$ws.Range("L25").Value = "198531.1 (±39891.6)"
$ws.Range("N25").Value = "1857.0 (±373.1)"

# Row 29
$ws.Range("C29").Value = 654.2
$ws.Range("E29").Value = 6944.5
$ws.Range("G29").Value = 10.62
$ws.Range("H29").Value = 1.03
$ws.Range("J29").Value = 993.3
$ws.Range("L29").Value = "6944.5 (±3326.7)"
$ws.Range("M29").Value = "10.62 (±1.03)"
$ws.Range("N29").Value = "993.3 (±475.8)"

# Row 30
$ws.Range("C30").Value = 496.9999999999999
$ws.Range("E30").Value = 5697.7
$ws.Range("G30").Value = 11.46
$ws.Range("H30").Value = 0.99
$ws.Range("J30").Value = 905.7
$ws.Range("L30").Value = "5697.7 (±2989.1)"
$ws.Range("M30").Value = "11.46 (±0.99)"
$ws.Range("N30").Value = "905.7 (±475.2)"

# Row 31
$ws.Range("C31").Value = 1109.6
$ws.Range("E31").Value = 12344.7
$ws.Range("G31").Value = 11.13
$ws.Range("H31").Value = 0.85
$ws.Range("J31").Value = 929.4
$ws.Range("L31").Value = "12344.7 (±5981.7)"
$ws.Range("M31").Value = "11.13 (±0.85)"
$ws.Range("N31").Value = "929.4 (±450.4)"

# Row 35
$ws.Range("C35").Value = 34442.39999999999
$ws.Range("E35").Value = 267999.1
$ws.Range("J35").Value = 771.4
$ws.Range("L35").Value = "267999.1 (±61937.7)"
$ws.Range("N35").Value = "771.4 (±178.3)"

# Row 36
$ws.Range("C36").Value = 35358.8
$ws.Range("E36").Value = 321601.0
$ws.Range("G36").Value = 9.1
$ws.Range("H36").Value = 0.36
$ws.Range("J36").Value = 989.5
$ws.Range("L36").Value = "321601.0 (±74176.6)"
$ws.Range("M36").Value = "9.1 (±0.36)"
$ws.Range("N36").Value = "989.5 (±228.2)"

# Row 37
$ws.Range("C37").Value = 69780.2
$ws.Range("E37").Value = 603257.8999999999
$ws.Range("G37").Value = 8.65
$ws.Range("H37").Value = 0.29
$ws.Range("J37").Value = 897.1
$ws.Range("L37").Value = "603257.9 (±136360.7)"
$ws.Range("M37").Value = "8.65 (±0.29)"
$ws.Range("N37").Value = "897.1 (±202.8)"

# Row 51
$ws.Range("C51").Value = 52652.8
$ws.Range("E51").Value = 412124.7
$ws.Range("J51").Value = 1414.9
$ws.Range("L51").Value = "412124.7 (±77971.2)"
$ws.Range("N51").Value = "1414.9 (±267.7)"

# Row 52
$ws.Range("C52").Value = 58262.0
$ws.Range("E52").Value = 576682.2
$ws.Range("G52").Value = 9.9
$ws.Range("H52").Value = 0.18
$ws.Range("J52").Value = 2035.9
$ws.Range("L52").Value = "576682.2 (±71607.0)"
$ws.Range("M52").Value = "9.9 (±0.18)"
$ws.Range("N52").Value = "2035.9 (±252.8)"

# Row 53
$ws.Range("C53").Value = 110910.6
$ws.Range("E53").Value = 1017432.2
$ws.Range("H53").Value = 0.22
$ws.Range("J53").Value = 1705.9
$ws.Range("L53").Value = "1017432.2 (±145489.8)"
$ws.Range("M53").Value = "9.17 (±0.22)"
$ws.Range("N53").Value = "1705.9 (±243.9)"

# Row 54
$ws.Range("C54").Value = 828.4000000000001
$ws.Range("D54").Value = 336.5
$ws.Range("E54").Value = 5535.9
$ws.Range("F54").Value = 2511.4
$ws.Range("G54").Value = 6.68
$ws.Range("H54").Value = 0.23
$ws.Range("J54").Value = 539.3
$ws.Range("K54").Value = 244.6
$ws.Range("L54").Value = "5535.9 (±2511.4)"
$ws.Range("M54").Value = "6.68 (±0.23)"
$ws.Range("N54").Value = "539.3 (±244.6)"

# Row 55
$ws.Range("C55").Value = 657.8
$ws.Range("D55").Value = 265.3
$ws.Range("E55").Value = 5439.5
$ws.Range("F55").Value = 2936.5
$ws.Range("G55").Value = 8.27
$ws.Range("H55").Value = 0.8
$ws.Range("J55").Value = 617.5
$ws.Range("K55").Value = 333.4
$ws.Range("L55").Value = "5439.5 (±2936.5)"
$ws.Range("M55").Value = "8.27 (±0.8)"
$ws.Range("N55").Value = "617.5 (±333.4)"

# Row 56
$ws.Range("C56").Value = 1450.6
$ws.Range("E56").Value = 10887.4
$ws.Range("G56").Value = 7.51
$ws.Range("H56").Value = 0.49
$ws.Range("J56").Value = 570.8
$ws.Range("L56").Value = "10887.4 (±5804.5)"
$ws.Range("M56").Value = "7.51 (±0.49)"
$ws.Range("N56").Value = "570.8 (±304.3)"

# Row 63
$ws.Range("C63").Value = 261.8
$ws.Range("E63").Value = 2499.8
$ws.Range("G63").Value = 9.55
$ws.Range("H63").Value = 1.92
$ws.Range("J63").Value = 1006.0
$ws.Range("L63").Value = "2499.8 (±1474.6)"
$ws.Range("M63").Value = "9.55 (±1.92)"
$ws.Range("N63").Value = "1006.0 (±593.4)"

# Row 64
$ws.Range("C64").Value = 300.2
$ws.Range("E64").Value = 4015.3
$ws.Range("G64").Value = 13.38
$ws.Range("H64").Value = 1.01
$ws.Range("J64").Value = 1512.1
$ws.Range("L64").Value = "4015.3 (±1816.9)"
$ws.Range("M64").Value = "13.38 (±1.01)"
$ws.Range("N64").Value = "1512.1 (±684.2)"

# Row 65
$ws.Range("C65").Value = 556.0
$ws.Range("E65").Value = 6477.3
$ws.Range("G65").Value = 11.65
$ws.Range("J65").Value = 1260.1
$ws.Range("L65").Value = "6477.3 (±2852.5)"
$ws.Range("M65").Value = "11.65 (±1.37)"
$ws.Range("N65").Value = "1260.1 (±554.9)"

# Row 66
$ws.Range("C66").Value = 337.0
$ws.Range("E66").Value = 3834.800000000001
$ws.Range("G66").Value = 11.38
$ws.Range("H66").Value = 0.8
$ws.Range("J66").Value = 1220.5
$ws.Range("L66").Value = "3834.8 (±1461.2)"
$ws.Range("M66").Value = "11.38 (±0.8)"
$ws.Range("N66").Value = "1220.5 (±465.0)"

# Row 67
$ws.Range("C67").Value = 543.0
$ws.Range("E67").Value = 5773.799999999999
$ws.Range("G67").Value = 10.63
$ws.Range("H67").Value = 0.83
$ws.Range("J67").Value = 1878.6
$ws.Range("L67").Value = "5773.8 (±2425.0)"
$ws.Range("M67").Value = "10.63 (±0.83)"
$ws.Range("N67").Value = "1878.6 (±789.0)"

# Row 68
$ws.Range("C68").Value = 859.0
$ws.Range("E68").Value = 9637.8
$ws.Range("G68").Value = 11.22
$ws.Range("H68").Value = 0.86
$ws.Range("J68").Value = 1550.6
$ws.Range("L68").Value = "9637.8 (±2992.0)"
$ws.Range("M68").Value = "11.22 (±0.86)"
$ws.Range("N68").Value = "1550.6 (±481.4)"

# Row 69
$ws.Range("C69").Value = 8504.2
$ws.Range("D69").Value = 1334.6
$ws.Range("E69").Value = 76458.0
$ws.Range("F69").Value = 19220.1
$ws.Range("H69").Value = 0.73
$ws.Range("J69").Value = 873.3
$ws.Range("K69").Value = 219.5
$ws.Range("L69").Value = "76458.0 (±19220.1)"
$ws.Range("M69").Value = "8.99 (±0.73)"
$ws.Range("N69").Value = "873.3 (±219.5)"

# Row 70
$ws.Range("C70").Value = 11933.6
$ws.Range("D70").Value = 1808.9
$ws.Range("E70").Value = 114090.3
$ws.Range("F70").Value = 22376.0
$ws.Range("G70").Value = 9.56
$ws.Range("H70").Value = 0.37
$ws.Range("J70").Value = 1320.1
$ws.Range("K70").Value = 258.9
$ws.Range("L70").Value = "114090.3 (±22376.0)"
$ws.Range("M70").Value = "9.56 (±0.37)"
$ws.Range("N70").Value = "1320.1 (±258.9)"

# Row 71
$ws.Range("C71").Value = 20342.6
$ws.Range("D71").Value = 3027.9
$ws.Range("E71").Value = 190631.7
$ws.Range("F71").Value = 41656.0
$ws.Range("G71").Value = 9.37
$ws.Range("H71").Value = 0.57
$ws.Range("J71").Value = 1095.7
$ws.Range("K71").Value = 239.4
$ws.Range("L71").Value = "190631.7 (±41656.0)"
$ws.Range("M71").Value = "9.37 (±0.57)"
$ws.Range("N71").Value = "1095.7 (±239.4)"

# Row 72
$ws.Range("C72").Value = 297.6
$ws.Range("D72").Value = 236.9
$ws.Range("E72").Value = 4613.7
$ws.Range("F72").Value = 3964.2
$ws.Range("G72").Value = 15.5
$ws.Range("H72").Value = 0.55
$ws.Range("J72").Value = 173.4
$ws.Range("K72").Value = 149.0
$ws.Range("L72").Value = "4613.7 (±3964.2)"
$ws.Range("M72").Value = "15.5 (±0.55)"
$ws.Range("N72").Value = "173.4 (±149.0)"

# Row 73
$ws.Range("C73").Value = 578.0
$ws.Range("E73").Value = 6853.5
$ws.Range("F73").Value = 3521.4
$ws.Range("G73").Value = 11.86
$ws.Range("H73").Value = 0.22
$ws.Range("J73").Value = 253.3
$ws.Range("L73").Value = "6853.5 (±3521.4)"
$ws.Range("M73").Value = "11.86 (±0.22)"
$ws.Range("N73").Value = "253.3 (±130.1)"

# Row 74
$ws.Range("C74").Value = 767.9999999999999
$ws.Range("D74").Value = 527.2
$ws.Range("E74").Value = 10676.7
$ws.Range("F74").Value = 7230.699999999999
$ws.Range("G74").Value = 13.9
$ws.Range("H74").Value = 0.07
$ws.Range("J74").Value = 198.9
$ws.Range("K74").Value = 134.7
$ws.Range("L74").Value = "10676.7 (±7230.7)"
$ws.Range("M74").Value = "13.9 (±0.07)"
$ws.Range("N74").Value = "198.9 (±134.7)"

# Row 81
$ws.Range("C81").Value = 16448.8
$ws.Range("E81").Value = 183274.3
$ws.Range("J81").Value = 1859.3
$ws.Range("L81").Value = "183274.3 (±36133.3)"
$ws.Range("N81").Value = "1859.3 (±366.6)"

# Row 82
$ws.Range("C82").Value = 23960.6
$ws.Range("E82").Value = 282223.9
$ws.Range("J82").Value = 2984.9
$ws.Range("L82").Value = "282223.9 (±44743.1)"
$ws.Range("N82").Value = "2984.9 (±473.2)"

# Row 83
$ws.Range("C83").Value = 40409.4
$ws.Range("E83").Value = 485440.7
$ws.Range("H83").Value = 0.01
$ws.Range("J83").Value = 2513.6
$ws.Range("L83").Value = "485440.7 (±79512.2)"
$ws.Range("M83").Value = "12.01 (±0.01)"
$ws.Range("N83").Value = "2513.6 (±411.7)"

# Row 87
$ws.Range("C87").Value = 3049.8
$ws.Range("E87").Value = 34360.2
$ws.Range("J87").Value = 1231.1
$ws.Range("L87").Value = "34360.2 (±10003.5)"
$ws.Range("N87").Value = "1231.1 (±358.4)"

# Row 88
$ws.Range("C88").Value = 3654.8
$ws.Range("E88").Value = 38994.4
$ws.Range("G88").Value = 10.67
$ws.Range("H88").Value = 0.43
$ws.Range("J88").Value = 1463.5
$ws.Range("L88").Value = "38994.4 (±7280.0)"
$ws.Range("M88").Value = "10.67 (±0.43)"
$ws.Range("N88").Value = "1463.5 (±273.2)"

# Row 89
$ws.Range("C89").Value = 6661.199999999999
$ws.Range("E89").Value = 74381.8
$ws.Range("G89").Value = 11.17
$ws.Range("H89").Value = 0.14
$ws.Range("J89").Value = 1363.4
$ws.Range("L89").Value = "74381.8 (±11753.5)"
$ws.Range("M89").Value = "11.17 (±0.14)"
$ws.Range("N89").Value = "1363.4 (±215.4)"

# Row 90
$ws.Range("C90").Value = 2183.8
$ws.Range("E90").Value = 16728.0
$ws.Range("G90").Value = 7.66
$ws.Range("H90").Value = 0.52
$ws.Range("J90").Value = 1602.5
$ws.Range("L90").Value = "16728.0 (±4050.0)"
$ws.Range("M90").Value = "7.66 (±0.52)"
$ws.Range("N90").Value = "1602.5 (±388.0)"

# Row 91
$ws.Range("C91").Value = 1950.8
$ws.Range("E91").Value = 18560.5
$ws.Range("G91").Value = 9.51
$ws.Range("H91").Value = 0.73
$ws.Range("J91").Value = 1767.4
$ws.Range("L91").Value = "18560.5 (±5002.2)"
$ws.Range("M91").Value = "9.51 (±0.73)"
$ws.Range("N91").Value = "1767.4 (±476.3)"

# Row 92
$ws.Range("C92").Value = 4122.4
$ws.Range("E92").Value = 35664.3
$ws.Range("G92").Value = 8.65
$ws.Range("J92").Value = 1703.2
$ws.Range("L92").Value = "35664.3 (±7026.4)"
$ws.Range("M92").Value = "8.65 (±0.41)"
$ws.Range("N92").Value = "1703.2 (±335.5)"

# Row 93
$ws.Range("C93").Value = 39333.0
$ws.Range("D93").Value = 4956.8
$ws.Range("E93").Value = 350345.8
$ws.Range("F93").Value = 54610.60000000001
$ws.Range("G93").Value = 8.91
$ws.Range("H93").Value = 0.23
$ws.Range("J93").Value = 1453.5
$ws.Range("K93").Value = 226.6
$ws.Range("L93").Value = "350345.8 (±54610.6)"
$ws.Range("M93").Value = "8.91 (±0.23)"
$ws.Range("N93").Value = "1453.5 (±226.6)"

# Row 94
$ws.Range("C94").Value = 39228.2
$ws.Range("E94").Value = 406580.8
$ws.Range("G94").Value = 10.36
$ws.Range("H94").Value = 0.14
$ws.Range("J94").Value = 1754.1
$ws.Range("L94").Value = "406580.8 (±59430.1)"
$ws.Range("M94").Value = "10.36 (±0.14)"
$ws.Range("N94").Value = "1754.1 (±256.4)"

# Row 95
$ws.Range("C95").Value = 78535.2
$ws.Range("D95").Value = 9621.2
$ws.Range("E95").Value = 772956.7999999999
$ws.Range("F95").Value = 110853.0
$ws.Range("G95").Value = 9.84
$ws.Range("H95").Value = 0.19
$ws.Range("J95").Value = 1634.7
$ws.Range("K95").Value = 234.4
$ws.Range("L95").Value = "772956.8 (±110853.0)"
$ws.Range("M95").Value = "9.84 (±0.19)"
$ws.Range("N95").Value = "1634.7 (±234.4)"
